$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: B2/D2 updated, C2/E2 cleared (removed)
$ws.Range("B2").Value = 2.0852145447608303
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2.2990046799629509
$ws.Range("E2").ClearContents()

# Row 3: all four values updated
$ws.Range("B3").Value = 1.7529298586155688
$ws.Range("C3").Value = -0.8956985217115051
$ws.Range("D3").Value = 1.6952949500695462
$ws.Range("E3").Value = -1.6748434028007984

# Update selection to match new used range for rows/cols B1:E3
$ws.Range("B1:E3").Select()
